$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) The document currently has a "_GoBack" bookmark sitting between the
#    runs " de manutenção" and ", citados na questão anterior." (last
#    paragraph). The edit relocates this bookmark to the end of the
#    title paragraph instead, so first remove it from its old spot.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) Title paragraph: "2ª Entrevista<trailing space>" -> "2ª Entrevista"
#    (underlined), followed by the relocated "_GoBack" bookmark.
#    We locate the trailing whitespace character just before the
#    paragraph mark, wrap a new bookmark around it, then delete that
#    character - leaving the (now empty) bookmark anchored immediately
#    after the run, exactly where the trailing space used to be.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1).Range
$trailStart = $titlePara.End - 2  # -1 for the paragraph mark, -1 more for the trailing space char
$trailRng = $d.Range($trailStart, $trailStart + 1)
if ($trailRng.Text -match '\s') {
    $d.Bookmarks.Add("_GoBack", $trailRng)
    $trailRng2 = $d.Range($trailStart, $trailStart + 1)
    $trailRng2.Delete()
}

# Underline the title paragraph (paragraph mark + run) as in the diff.
$titlePara2 = $d.Paragraphs.Item(1).Range
$titlePara2.Font.Underline = 1

# ---------------------------------------------------------------------
# 3) Question 8: merge the "8. " run with the following run into a
#    single run. A Find/Replace that spans both runs causes the engine
#    to coalesce them (they already share identical formatting).
# ---------------------------------------------------------------------
$d.Content.Find.Execute( `
    "8. Após a compra do veículo, o que acontece até o veículo ser disponibilizado para compra?", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "8. Após a compra do veículo, o que acontece até o veículo ser disponibilizado para compra?", `
    2) | Out-Null

# ---------------------------------------------------------------------
# 4) Last paragraph: merge " de manutenção" with ", citados na questão
#    anterior." into a single run, now that the bookmark that used to
#    separate them has been removed. Restricting the Find text to the
#    second run only (it no longer has a barrier after it) keeps the
#    preceding "Antes do veículo...serviços" run untouched.
# ---------------------------------------------------------------------
$d.Content.Find.Execute( `
    ", citados na questão anterior.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    ", citados na questão anterior.", `
    2) | Out-Null
